$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value for the Price (D) and Volume(1h) (E) columns.
# Force each touched cell to Text format first so the numeric-looking /
# percent-looking strings are not auto-converted to Number/Percent by Excel,
# matching the original inlineStr (text) cell content.
$updates = @(
    @{ Cell = "D2"; Value = "325.93" }
    @{ Cell = "E2"; Value = "-2.78%" }
    @{ Cell = "D3"; Value = "44.52" }
    @{ Cell = "E3"; Value = "1.49%" }
    @{ Cell = "D4"; Value = "5.577" }
    @{ Cell = "E4"; Value = "-4.08%" }
    @{ Cell = "D5"; Value = "0.08077" }
    @{ Cell = "E5"; Value = "-3.14%" }
    @{ Cell = "D6"; Value = "8.678" }
    @{ Cell = "D7"; Value = "1.905" }
    @{ Cell = "E7"; Value = "-4.42%" }
    @{ Cell = "D8"; Value = "4.293" }
    @{ Cell = "E8"; Value = "-4.67%" }
    @{ Cell = "D9"; Value = "2.690" }
    @{ Cell = "E9"; Value = "-7.18%" }
    @{ Cell = "D10"; Value = "0.9425" }
    @{ Cell = "E10"; Value = "0.31%" }
    @{ Cell = "D11"; Value = "0.1172" }
    @{ Cell = "E11"; Value = "-5.40%" }
    @{ Cell = "D12"; Value = "0.1857" }
    @{ Cell = "E12"; Value = "-4.74%" }
    @{ Cell = "D13"; Value = "0.1000" }
    @{ Cell = "E13"; Value = "3.44%" }
    @{ Cell = "D14"; Value = "0.04278" }
    @{ Cell = "E14"; Value = "-6.66%" }
    @{ Cell = "E15"; Value = "-0.26%" }
    @{ Cell = "D16"; Value = "0.001280" }
    @{ Cell = "E16"; Value = "-1.31%" }
    @{ Cell = "D17"; Value = "0.04193" }
    @{ Cell = "E17"; Value = "-4.68%" }
    @{ Cell = "D18"; Value = "0.005881" }
    @{ Cell = "E18"; Value = "-2.48%" }
    @{ Cell = "D19"; Value = "3.576" }
    @{ Cell = "E19"; Value = "2.29%" }
    @{ Cell = "D20"; Value = "0.3499" }
    @{ Cell = "E20"; Value = "-0.31%" }
    @{ Cell = "D21"; Value = "8.443" }
    @{ Cell = "E21"; Value = "-3.69%" }
    @{ Cell = "D22"; Value = "0.1369" }
    @{ Cell = "E22"; Value = "0.54%" }
    @{ Cell = "D24"; Value = "0.001242" }
    @{ Cell = "E24"; Value = "-1.11%" }
    @{ Cell = "D25"; Value = "0.004510" }
    @{ Cell = "E25"; Value = "2.69%" }
    @{ Cell = "E26"; Value = "-6.26%" }
    @{ Cell = "D27"; Value = "0.0003989" }
    @{ Cell = "E27"; Value = "-0.05%" }
    @{ Cell = "D39"; Value = "0.02639" }
    @{ Cell = "E39"; Value = "-5.97%" }
    @{ Cell = "D40"; Value = "0.05444" }
    @{ Cell = "E40"; Value = "-4.85%" }
    @{ Cell = "D41"; Value = "0.007654" }
    @{ Cell = "E41"; Value = "-3.37%" }
    @{ Cell = "D42"; Value = "0.1398" }
    @{ Cell = "E42"; Value = "-2.21%" }
    @{ Cell = "D43"; Value = "0.007070" }
    @{ Cell = "E43"; Value = "-21.16%" }
    @{ Cell = "D44"; Value = "0.002026" }
    @{ Cell = "E44"; Value = "-6.03%" }
    @{ Cell = "D45"; Value = "0.008851" }
    @{ Cell = "E45"; Value = "-16.08%" }
    @{ Cell = "D46"; Value = "0.00007168" }
    @{ Cell = "E46"; Value = "-0.58%" }
    @{ Cell = "E47"; Value = "0.10%" }
    @{ Cell = "D48"; Value = "0.003665" }
    @{ Cell = "E48"; Value = "13.02%" }
    @{ Cell = "D49"; Value = "0.002270" }
    @{ Cell = "E49"; Value = "-0.45%" }
    @{ Cell = "D50"; Value = "0.00002103" }
    @{ Cell = "E50"; Value = "0.10%" }
    @{ Cell = "E51"; Value = "0.10%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

Write-Host "Updated" $updates.Count "cells"
